# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# --- Rushing sheet updates ---
# Row 2: P.Mahomes
$rushing.Range("C2").Value = 16
$rushing.Range("D2").Value = 18
$rushing.Range("E2").Value = 9
$rushing.Range("F2").Value = 13

# Row 5: Da.Williams
$rushing.Range("C5").Value = 52
$rushing.Range("D5").Value = 26
$rushing.Range("E5").Value = 10

# Row 6: J.McKinnon
$rushing.Range("C6").Value = 16
$rushing.Range("D6").Value = 4
$rushing.Range("F6").Value = 3

# Row 7: D.Gore
$rushing.Range("C7").Value = 26
$rushing.Range("F7").Value = 4

# Row 8: M.Burton
$rushing.Range("E8").Value = 6
$rushing.Range("F8").Value = 4

# Row 10: M.Hardman
$rushing.Range("C10").Value = 8

# --- Receiving sheet updates ---
# Row 3: Da.Williams
$receiving.Range("C3").Value = 34
$receiving.Range("D3").Value = 25
$receiving.Range("G3").Value = 6

# Row 4: J.McKinnon
$receiving.Range("C4").Value = 9
$receiving.Range("D4").Value = 8
$receiving.Range("E4").Value = 3
$receiving.Range("F4").Value = 2
$receiving.Range("G4").Value = 4
$receiving.Range("H4").Value = 3

# Row 5: D.Gore
$receiving.Range("C5").Value = 5
$receiving.Range("G5").Value = 1

# Row 7: T.Hill
$receiving.Range("C7").Value = 124
$receiving.Range("D7").Value = 98
$receiving.Range("E7").Value = 37
$receiving.Range("F7").Value = 15
$receiving.Range("G7").Value = 27
$receiving.Range("H7").Value = 18

# Row 8: M.Hardman
$receiving.Range("C8").Value = 64
$receiving.Range("D8").Value = 50
$receiving.Range("E8").Value = 18
$receiving.Range("F8").Value = 9
$receiving.Range("G8").Value = 14
$receiving.Range("H8").Value = 7

# Row 9: B.Pringle
$receiving.Range("C9").Value = 44
$receiving.Range("D9").Value = 30
$receiving.Range("E9").Value = 15
$receiving.Range("G9").Value = 5
$receiving.Range("H9").Value = 3

# Row 10: D.Robinson
$receiving.Range("C10").Value = 31
$receiving.Range("D10").Value = 21
$receiving.Range("E10").Value = 11
$receiving.Range("F10").Value = 5

# Row 12: J.Gordon
$receiving.Range("C12").Value = 8

# Row 13: T.Kelce
$receiving.Range("C13").Value = 117
$receiving.Range("D13").Value = 83
$receiving.Range("E13").Value = 25
$receiving.Range("F13").Value = 15
$receiving.Range("G13").Value = 20
$receiving.Range("H13").Value = 16

# Row 14: B.Bell
$receiving.Range("C14").Value = 12
$receiving.Range("D14").Value = 8
$receiving.Range("G14").Value = 2
$receiving.Range("H14").Value = 2

# Row 15: N.Gray
$receiving.Range("C15").Value = 9
